# Adds a tooltip row for "Wind Speed" right after the existing "wind" row,
# matching the commit: "added tooltip for humidity and wind speed icon for clarity"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 44 (pushes winter/summer/etc. down by one row)
$ws.Rows.Item(44).Insert()

# Populate the new row with the windSpeed key/value pair
$ws.Range("A44").Value = "windSpeed"
$ws.Range("B44").Value = "Wind Speed"

# Reflect the cursor / selection state left behind by the edit
[void]$ws.Application.Goto($ws.Range("A31"), $true)
[void]$ws.Range("B44").Select()
